# Calibrate ETS for split between onshore and offshore wind
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETS")

# Onshore wind (row 6): 0.7 -> 0.65 across all year columns (B:AF)
$ws.Range("B6:AF6").Value = 0.65

# Solar pv (row 7): 0.5 -> 0.55 across all year columns (B:AF)
$ws.Range("B7:AF7").Value = 0.55

# Offshore wind (row 14): 0.3 -> 1 across all year columns (B:AF)
$ws.Range("B14:AF14").Value = 1

# Reflect the user's selection at save time (B6:AF6), matching the
# updated sheetView selection/topLeftCell recorded in the workbook.
$ws.Range("B6:AF6").Select()
